$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-12 (header + 11 data rows) are unchanged.

# New rows appended to the data dictionary table (row 13 used to hold the
# empty "end of table" marker cell - it now becomes a real data row, and
# the marker moves down to row 17).
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "TotalMinsIntnl"
$ws.Range("C13").Value = "Número total de minutos en llamadas internacionales"
$ws.Range("B13").Style = "Normal"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "TotalLlamadasIntnl"
$ws.Range("C14").Value = "Número total de llamadas internacionales"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "TotalLlamadasCallCenter"
$ws.Range("C15").Value = "Número total de llamadas al call center"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Estado"
$ws.Range("C16").Value = " Indica si el cliente está vinculado o retirado"

# The empty bold marker cell moves down to just past the new last row.
$ws.Range("B17").Font.Bold = $true

$excel.Goto($ws.Range("A2:C16"))
